$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before the current row 260, pushing the existing
# rows 260-265 down to 262-267.
$ws.Rows("260:261").Insert()

# New row 260 data
$ws.Cells.Item(260, 1).Value = 10
$ws.Cells.Item(260, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(260, 3).Value = "La Araucanía"
$ws.Cells.Item(260, 4).Value = 44448
$ws.Cells.Item(260, 5).Value = 9
$ws.Cells.Item(260, 6).Value = "Fruta"
$ws.Cells.Item(260, 7).Value = 100108
$ws.Cells.Item(260, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(260, 9).Value = 100108005
$ws.Cells.Item(260, 10).Value = "Piña"
$ws.Cells.Item(260, 11).Value = "Caramelo"
$ws.Cells.Item(260, 12).Value = "Primera"
$ws.Cells.Item(260, 13).Value = 100
$ws.Cells.Item(260, 14).Value = 20000
$ws.Cells.Item(260, 15).Value = 20000
$ws.Cells.Item(260, 16).Value = 20000
$ws.Cells.Item(260, 17).Value = "`$/caja 12 unidades"
$ws.Cells.Item(260, 18).Value = "Ecuador"
$ws.Cells.Item(260, 19).Value = 1667
$ws.Cells.Item(260, 20).Value = 12

# New row 261 data
$ws.Cells.Item(261, 1).Value = 10
$ws.Cells.Item(261, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(261, 3).Value = "La Araucanía"
$ws.Cells.Item(261, 4).Value = 44448
$ws.Cells.Item(261, 5).Value = 9
$ws.Cells.Item(261, 6).Value = "Fruta"
$ws.Cells.Item(261, 7).Value = 100108
$ws.Cells.Item(261, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(261, 9).Value = 100108005
$ws.Cells.Item(261, 10).Value = "Piña"
$ws.Cells.Item(261, 11).Value = "Caramelo"
$ws.Cells.Item(261, 12).Value = "Segunda"
$ws.Cells.Item(261, 13).Value = 70
$ws.Cells.Item(261, 14).Value = 20000
$ws.Cells.Item(261, 15).Value = 21000
$ws.Cells.Item(261, 16).Value = 20571
$ws.Cells.Item(261, 17).Value = "`$/caja 14 unidades"
$ws.Cells.Item(261, 18).Value = "Ecuador"
$ws.Cells.Item(261, 19).Value = 1469
$ws.Cells.Item(261, 20).Value = 14
